$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "RA40"
$ws.Range("B3").Value = "Ben Suggs"
$ws.Range("C3").Value = "Happy that Cameron is here to deal with this. "
$ws.Range("D3").Value = "2025-09-30 19:45:31"
